$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAC-10")

# Update the synonym text for "Sucursal" (row 12, column C) from "Sede" to "Sede, oficina"
$ws.Range("C12").Value = "Sede, oficina"

# Move the active selection from B7 to F9
$ws.Activate()
$ws.Range("F9").Select()
